# Update the "Prix Spot" sheet: append a new day column (BV) of data,
# mirroring the existing 25-aug (BU) column's header/data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell BV1: copy the formatting (style) of BU1 so the new header
# reuses the same existing header style, then set its new date label.
$ws.Range("BU1").Copy()
$ws.Range("BV1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BV1").Value = "26-aug"

# New daily price values for rows 2-25 (hours 00-01 .. 23-24)
$values = @(
    94.83,
    86.78,
    70.37,
    63.58,
    58.35,
    66.86,
    86.97,
    95.44,
    97.48,
    90,
    70,
    60.65,
    41.25,
    26.28,
    38.05,
    47.44,
    60,
    65.03,
    76.09999999999999,
    103.14,
    117,
    123.46,
    109.5,
    97.90000000000001
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 74).Value = $values[$i]
}
